$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 126.666664
$ws.Range("I6").Value = 126.666664
$ws.Range("K6").Value = 379.999992
$ws.Range("M6").Value = -267.999992
$ws.Range("H51").Value = 6272
$ws.Range("I51").Value = 9249
$ws.Range("J51").Value = 5610.4443
$ws.Range("K51").Value = 9249
$ws.Range("L51").Value = 5610.4443
$ws.Range("M51").Value = -8765
$ws.Range("N51").Value = -6578.4443
$ws.Range("H132").Value = 2096.25
$ws.Range("I132").Value = 2096.25
$ws.Range("K132").Value = 6288.75
$ws.Range("M132").Value = -3758.75

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H53").Value = 20000
$ws.Range("I53").Value = 20000
$ws.Range("K53").Value = 20000
$ws.Range("M53").Value = -19318
$ws.Range("H97").Value = 1417.3334
$ws.Range("I97").Value = 1126
$ws.Range("J97").Value = 2000
$ws.Range("K97").Value = 1126
$ws.Range("L97").Value = 2000
$ws.Range("M97").Value = -630
$ws.Range("N97").Value = -2992
$ws.Range("H112").Value = 27789.5
$ws.Range("J112").Value = 27789.5
$ws.Range("L112").Value = 27789.5
$ws.Range("N112").Value = -30743.5

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 235.22223
$ws.Range("I22").Value = 235.22223
$ws.Range("K22").Value = 235.22223
$ws.Range("M22").Value = -62.22223
$ws.Range("H26").Value = 20000
$ws.Range("I26").Value = 20000
$ws.Range("K26").Value = 20000
$ws.Range("M26").Value = -19708
$ws.Range("H86").Value = 1771.3334
$ws.Range("I86").Value = 1771.3334
$ws.Range("K86").Value = 1771.3334
$ws.Range("M86").Value = -648.3334
$ws.Range("H89").Value = 1771.3334
$ws.Range("I89").Value = 1771.3334
$ws.Range("K89").Value = 8856.666999999999
$ws.Range("M89").Value = -3240.666999999999
$ws.Range("H96").Value = 20000
$ws.Range("I96").Value = 20000
$ws.Range("K96").Value = 20000
$ws.Range("M96").Value = -17254
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2465.4285
$ws.Range("I31").Value = 1851.8
$ws.Range("K31").Value = 1851.8
$ws.Range("M31").Value = -1556.8
$ws.Range("H34").Value = 2465.4285
$ws.Range("I34").Value = 1851.8
$ws.Range("K34").Value = 1851.8
$ws.Range("M34").Value = -1649.8
$ws.Range("H36").Value = 4000
$ws.Range("J36").Value = 4000
$ws.Range("L36").Value = 4000
$ws.Range("N36").Value = -4776
$ws.Range("H40").Value = 4000
$ws.Range("J40").Value = 4000
$ws.Range("L40").Value = 4000
$ws.Range("N40").Value = -4320
$ws.Range("H42").Value = 26750
$ws.Range("I42").Value = 26750
$ws.Range("K42").Value = 26750
$ws.Range("M42").Value = -26157
$ws.Range("H44").Value = 2933.3333
$ws.Range("I44").Value = 1900
$ws.Range("J44").Value = 5000
$ws.Range("K44").Value = 1900
$ws.Range("L44").Value = 5000
$ws.Range("M44").Value = -1458
$ws.Range("N44").Value = -5884
$ws.Range("H55").Value = 21666.666
$ws.Range("I55").Value = 10000
$ws.Range("J55").Value = 27500
$ws.Range("K55").Value = 10000
$ws.Range("L55").Value = 27500
$ws.Range("M55").Value = -9685
$ws.Range("N55").Value = -28130
$ws.Range("H103").Value = 15000
$ws.Range("I103").Value = 15000
$ws.Range("K103").Value = 15000
$ws.Range("M103").Value = -13828
$ws.Range("H105").Value = 811.5
$ws.Range("I105").Value = 875
$ws.Range("J105").Value = 663.3333
$ws.Range("K105").Value = 875
$ws.Range("L105").Value = 663.3333
$ws.Range("M105").Value = 872
$ws.Range("N105").Value = -4157.3333
$ws.Range("H132").Value = 5285.7144
$ws.Range("I132").Value = 4746.5
$ws.Range("K132").Value = 14239.5
$ws.Range("M132").Value = -11709.5
$ws.Range("H134").Value = 8166.5
$ws.Range("J134").Value = 8166.5
$ws.Range("L134").Value = 24499.5
$ws.Range("N134").Value = -29569.5

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("M11").ClearContents()
$ws.Range("H92").Value = 134.4
$ws.Range("J92").Value = 105.5
$ws.Range("L92").Value = 316.5
$ws.Range("N92").Value = -2812.5
$ws.Range("H113").Value = 433.75
$ws.Range("J113").Value = 295
$ws.Range("L113").Value = 885
$ws.Range("N113").Value = -5225

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 10000
$ws.Range("J19").Value = 10000
$ws.Range("L19").Value = 10000
$ws.Range("N19").Value = -10576

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()
$ws.Range("H22").Value = 999.5
$ws.Range("I22").Value = 999
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 999
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -704
$ws.Range("N22").Value = -1590
$ws.Range("H27").Value = 999.5
$ws.Range("I27").Value = 999
$ws.Range("J27").Value = 1000
$ws.Range("K27").Value = 999
$ws.Range("L27").Value = 1000
$ws.Range("M27").Value = -892
$ws.Range("N27").Value = -1214
$ws.Range("H45").Value = 3000
$ws.Range("I45").Value = 3000
$ws.Range("K45").Value = 3000
$ws.Range("M45").Value = -2593
$ws.Range("H55").Value = 2136.3333
$ws.Range("I55").Value = 378.42856
$ws.Range("J55").Value = 4597.4
$ws.Range("K55").Value = 378.42856
$ws.Range("L55").Value = 4597.4
$ws.Range("M55").Value = -205.42856
$ws.Range("N55").Value = -4943.4
$ws.Range("H110").Value = 26321.75
$ws.Range("J110").Value = 26321.75
$ws.Range("L110").Value = 26321.75
$ws.Range("N110").Value = -34501.75
$ws.Range("H122").Value = 5882.846
$ws.Range("I122").Value = 6047.8
$ws.Range("J122").Value = 5333
$ws.Range("K122").Value = 18143.4
$ws.Range("L122").Value = 15999
$ws.Range("M122").Value = -15693.4
$ws.Range("N122").Value = -20899
$ws.Range("H132").Value = 593
$ws.Range("I132").Value = 593
$ws.Range("K132").Value = 1779
$ws.Range("M132").Value = 751

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1754.9
$ws.Range("I132").Value = 1505.4445
$ws.Range("K132").Value = 4516.333500000001
$ws.Range("M132").Value = -1986.333500000001
